$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.956.47'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '3.423.96'
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.98'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.24'
$ws.Range("E6").Value = '  -2.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.631'
$ws.Range("E7").Value = '  +7.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.731'
$ws.Range("E9").Value = '  +7.10%  '
$ws.Range("E10").Value = '  +10.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.44'
$ws.Range("E11").Value = '  +1.76%  '
$ws.Range("E12").Value = '  +0.09%  '
$ws.Range("D13").Value = '3.959.05'
$ws.Range("E13").Value = '  +0.37%  '
$ws.Range("E14").Value = '  +7.85%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.12'
$ws.Range("E15").Value = '  +7.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000202'
$ws.Range("E16").Value = '  +44.11%  '
$ws.Range("D17").Value = '3.432.94'
$ws.Range("E17").Value = '  +0.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.34'
$ws.Range("E18").Value = '  +6.29%  '
$ws.Range("E19").Value = '  +6.53%  '
$ws.Range("D20").Value = '61.863.81'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '442.05'
$ws.Range("E21").Value = '  +42.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.15'
$ws.Range("E22").Value = '  +9.85%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.19'
$ws.Range("E23").Value = '  +1.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.87'
$ws.Range("E24").Value = '  +1.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +3.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '32.77'
$ws.Range("E26").Value = '  +11.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.73'
$ws.Range("E27").Value = '  +9.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.80'
$ws.Range("E28").Value = '  -0.55%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.76'
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.54'
$ws.Range("E30").Value = '  -6.60%  '
$ws.Range("E31").Value = '  +5.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.171'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("E33").Value = '  -0.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.67'
$ws.Range("E34").Value = '  -0.50%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0496'
$ws.Range("E36").Value = '  +3.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.10'
$ws.Range("E37").Value = '  +3.39%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("E40").Value = '  +7.31%  '
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("E42").Value = '  -2.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.09'
$ws.Range("E43").Value = '  +1.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.26'
$ws.Range("E44").Value = '  +8.10%  '
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("E46").Value = '  +8.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.49'
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.20'
$ws.Range("E48").Value = '  +4.74%  '
$ws.Range("D49").Value = '3.767.00'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.117.45'
$ws.Range("E50").Value = '  +0.99%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.06'
$ws.Range("E51").Value = '  +6.85%  '
